# "adding new progress as of date 04 nov 2025"
#
# On the "Training Dashboard" sheet, for every data row (3-37):
#   - column H (PERIOD TO EXPIRE) ticks down by one day
#   - column I (LAST UPDATE) moves from 03-Nov-2025 to 04-Nov-2025
#
# Column I stores its date as literal text (not a real Excel date), so a
# plain .Value assignment of "04-Nov-2025" would get auto-converted into a
# date serial by Excel's input parsing and would also stamp a new date
# number-format onto the cell. To keep the cell a plain text value (and its
# original style untouched) the literal is written through a formula that
# evaluates to a text string, then converted in place to a static value via
# Copy / PasteSpecial(Values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$newPeriods = @(626,483,483,285,356,467,489,470,472,664,397,300,529,278,311,377,378,453,-2,324,326,325,133,151,133,156,150,266,255,264,264,264,299,299,603)

$startRow = 3
$endRow = 37

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = $startRow + $i

    # Column H: PERIOD TO EXPIRE -> one day closer to expiry
    $ws.Cells.Item($row, 8).Value = $newPeriods[$i]

    # Column I: LAST UPDATE -> stage the new date as a text formula so it
    # isn't auto-parsed into a real date value.
    $ws.Cells.Item($row, 9).Formula = '="04-Nov-2025"'
}

# Collapse the staged formulas down to plain text values in one shot,
# preserving each cell's existing style/number format.
$rngI = $ws.Range("I$startRow`:I$endRow")
$rngI.Copy()
$rngI.PasteSpecial(-4163)
